$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.786.59"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "'3.751.19"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'622.08"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'180.23"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "'3.748.14"
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  +3.28%  "
$ws.Range("E11").Value = "  -5.10%  "
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").Value = "'40.92"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").Value = "'4.367.11"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "'3.754.69"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("D17").Value = "'69.803.48"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'7.64"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").Value = "'16.74"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").Value = "'505.86"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").Value = "'9.50"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D23").Value = "'0.730"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "'2.50"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").Value = "'13.17"
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("D27").Value = "'11.11"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  +26.00%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").Value = "'7.89"
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("D33").Value = "'31.18"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "'0.116"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("D37").Value = "'6.20"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "'0.336"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("D41").Value = "'49.95"
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("D42").Value = "'45.54"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'425.75"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'8.73"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("D46").Value = "'3.004.50"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'27.38"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'137.69"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = "  +1.47%  "
